# "Probando el codigo (Version 1)" -- add a second question (row 5) to the
# Preguntas sheet, fill in the two missing options for the first question
# (row 4), then leave the selection where the author left off (B13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 4 ("Se pueden crear Activity con XML?") with its remaining
# two options.
$ws.Range("D4").Value = "TODAS"
$ws.Range("E4").Value = "NINGUNA"

# New question in row 5.
$ws.Range("A5").Value = "¿Sere capaz de usar un recycler view adecuadmente?"
$ws.Range("B5").Value = "NO"
$ws.Range("C5").Value = "SI "
$ws.Range("D5").Value = "NO SE"
$ws.Range("E5").Value = "MAS O MENOS"

# Match where the cursor ended up when the author saved the file.
$ws.Range("B13").Select() | Out-Null
